$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $r = $ws.Range($cellRef)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.Style = "Normal"
}

# Row 2
$ws.Range('D2').Value = '72.677.79'
$ws.Range('E2').Value = '  +3.78%  '

# Row 3
$ws.Range('D3').Value = '2.635.71'
$ws.Range('E3').Value = '  +2.45%  '

# Row 4
$ws.Range('E4').Value = '  +0.16%  '

# Row 5
Set-TextValue 'D5' '605.50'
$ws.Range('E5').Value = '  +0.56%  '

# Row 6
Set-TextValue 'D6' '179.61'
$ws.Range('E6').Value = '  +0.42%  '

# Row 7
$ws.Range('E7').Value = '  +0.07%  '

# Row 8
Set-TextValue 'D8' '0.528'
$ws.Range('E8').Value = '  +1.43%  '

# Row 9
$ws.Range('E9').Value = '  +9.91%  '

# Row 10
$ws.Range('D10').Value = '2.636.93'
$ws.Range('E10').Value = '  +2.60%  '

# Row 11
$ws.Range('E11').Value = '  +1.08%  '

# Row 12
Set-TextValue 'D12' '0.356'
$ws.Range('E12').Value = '  +3.23%  '

# Row 13
Set-TextValue 'D13' '5.04'
$ws.Range('E13').Value = '  +0.41%  '

# Row 14
$ws.Range('E14').Value = '  +4.34%  '

# Row 15
$ws.Range('D15').Value = '3.108.02'
$ws.Range('E15').Value = '  +2.16%  '

# Row 16
$ws.Range('D16').Value = '72.635.62'
$ws.Range('E16').Value = '  +3.85%  '

# Row 17
Set-TextValue 'D17' '26.80'
$ws.Range('E17').Value = '  +1.58%  '

# Row 18
$ws.Range('D18').Value = '2.639.96'
$ws.Range('E18').Value = '  +2.27%  '

# Row 19
$ws.Range('B19').Value = 'Chainlink'
$ws.Range('C19').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
Set-TextValue 'D19' '11.76'
$ws.Range('E19').Value = '  +4.99%  '

# Row 20
$ws.Range('B20').Value = 'BitcoinCash'
$ws.Range('C20').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
Set-TextValue 'D20' '383.96'
$ws.Range('E20').Value = '  +4.47%  '

# Row 21
Set-TextValue 'D21' '7.95'
$ws.Range('E21').Value = '  +3.24%  '

# Row 22
$ws.Range('E22').Value = '  +2.32%  '

# Row 23
Set-TextValue 'D23' '2.06'
$ws.Range('E23').Value = '  +16.37%  '

# Row 24
Set-TextValue 'D24' '74.11'
$ws.Range('E24').Value = '  +4.34%  '

# Row 25
Set-TextValue 'D25' '4.42'
$ws.Range('E25').Value = '  +2.46%  '

# Row 26
$ws.Range('E26').Value = '  -0.02%  '

# Row 27
Set-TextValue 'D27' '10.04'
$ws.Range('E27').Value = '  +8.59%  '

# Row 28
$ws.Range('D28').Value = '2.773.16'
$ws.Range('E28').Value = '  +2.55%  '

# Row 29
$ws.Range('E29').Value = '  +0.07%  '

# Row 30
$ws.Range('D30').Value = '0.0₃0960'
$ws.Range('E30').Value = '  +3.52%  '

# Row 31
$ws.Range('E31').Value = '  +4.02%  '

# Row 32
Set-TextValue 'D32' '520.55'
$ws.Range('E32').Value = '  +0.12%  '

# Row 33
Set-TextValue 'D33' '1.33'
$ws.Range('E33').Value = '  +3.72%  '

# Row 34
Set-TextValue 'D34' '1.83'
$ws.Range('E34').Value = '  +1.32%  '

# Row 35
Set-TextValue 'D35' '1.00'
$ws.Range('E35').Value = '  +0.10%  '

# Row 36
Set-TextValue 'D36' '163.26'
$ws.Range('E36').Value = '  +0.05%  '

# Row 37
Set-TextValue 'D37' '19.44'
$ws.Range('E37').Value = '  +1.73%  '

# Row 38
Set-TextValue 'D38' '1.41'
$ws.Range('E38').Value = '  +4.07%  '

# Row 39
$ws.Range('B39').Value = 'WhiteBITCoin'
$ws.Range('C39').Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
Set-TextValue 'D39' '19.11'
$ws.Range('E39').Value = '  +0.97%  '

# Row 40
$ws.Range('B40').Value = 'Kaspa'
$ws.Range('C40').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
Set-TextValue 'D40' '0.112'
$ws.Range('E40').Value = '  -6.11%  '

# Row 41
Set-TextValue 'D41' '1.87'
$ws.Range('E41').Value = '  +5.53%  '

# Row 42
Set-TextValue 'D42' '5.18'
$ws.Range('E42').Value = '  +4.39%  '

# Row 43
$ws.Range('E43').Value = '  +0.04%  '

# Row 44
Set-TextValue 'D44' '2.58'
$ws.Range('E44').Value = '  +3.87%  '

# Row 45
Set-TextValue 'D45' '0.335'
$ws.Range('E45').Value = '  +2.37%  '

# Row 46
Set-TextValue 'D46' '39.45'
$ws.Range('E46').Value = '  +0.99%  '

# Row 47
Set-TextValue 'D47' '150.90'
$ws.Range('E47').Value = '  -1.85%  '

# Row 48
Set-TextValue 'D48' '3.71'
$ws.Range('E48').Value = '  +1.84%  '

# Row 49
Set-TextValue 'D49' '0.546'
$ws.Range('E49').Value = '  +3.92%  '

# Row 50
$ws.Range('E50').Value = '  +4.05%  '

# Row 51
$ws.Range('D51').Value = '0.0₆0265'
$ws.Range('E51').Value = '  +2.20%  '
